$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LinearRegression - only B2 changes
$ws.Range("B2").Value = 4618500226688193

# Row 3: RandomForestRegressor - B3, C3, D3 change
$ws.Range("B3").Value = 97459378798436.83
$ws.Range("C3").Value = 87641103758085.58
$ws.Range("D3").Value = 457284915657149.5

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 1808331252240.809
$ws.Range("C4").Value = 1776725103481.106
$ws.Range("D4").Value = 111761243991670.3

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 286363457254420.1
$ws.Range("C5").Value = 318283505370468.3
$ws.Range("D5").Value = 2534286757029580
